# Adding base-2 exponential demo.
#
# The workbook currently has two sheets: "Exponential" (f(x) = EXP(x)) and
# "Logarithm" (f(x) = LN(x)). This adds a new "Base 2 Exponential" sheet,
# positioned before "Exponential", that mirrors the "Exponential" sheet's
# layout/x-values but uses f(x) = 2^x instead of f(x) = EXP(x).

$wb = $excel.ActiveWorkbook

$expSheet = $wb.Worksheets.Item("Exponential")

# Duplicate the "Exponential" sheet (keeps headers, x-values, number
# formats/styles, column widths, page setup, etc. identical) and place the
# copy immediately before "Exponential" so tab order becomes:
#   Base 2 Exponential, Exponential, Logarithm
$expSheet.Copy($expSheet)

$newSheet = $wb.Worksheets.Item("Exponential (2)")
$newSheet.Name = "Base 2 Exponential"

# Re-point every f(x) formula in column B from EXP(Ax) to 2^Ax.
for ($row = 2; $row -le 82; $row++) {
    $newSheet.Range("B$row").Formula = "=2^A$row"
}
